$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.100.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4692"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8794"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.885.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07357"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.366"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.520"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008705"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.543.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.233"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.088.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.879"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.166"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08938"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7413"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.506"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.948"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.008"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.089"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05287"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.431"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.970"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.304"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1638"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.389"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4867"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06269"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.52%  "
